$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price groups were reported, each inserted as a 2-row block.
# Group "2022-07-06" goes at the very top of this product's data (row 156),
# pushing the existing rows down by 2.
$ws.Rows.Item(156).Insert()
$ws.Rows.Item(156).Insert()

# Group "2022-07-05" goes further down, right after the old row 163 data
# (which is now at row 165 after the first insertion above), i.e. at row 166.
$ws.Rows.Item(166).Insert()
$ws.Rows.Item(166).Insert()

# --- Row 156 ---
$ws.Range("A156").Value = 11
$ws.Range("B156").Value = 'Vega Monumental Concepción'
$ws.Range("C156").Value = 'Bíobío'
$ws.Range("D156").Value = 44748
$ws.Range("E156").Value = 8
$ws.Range("F156").Value = 'Fruta'
$ws.Range("G156").Value = 100101
$ws.Range("H156").Value = 'Berries'
$ws.Range("I156").Value = 100101007
$ws.Range("J156").Value = 'Kiwi'
$ws.Range("K156").Value = 'Hayward'
$ws.Range("L156").Value = 'Primera'
$ws.Range("M156").Value = 50
$ws.Range("N156").Value = 8000
$ws.Range("O156").Value = 8000
$ws.Range("P156").Value = 8000
$ws.Range("Q156").Value = '$/bandeja 18 kilos'
$ws.Range("R156").Value = 'Región de O''Higgins'
$ws.Range("S156").Value = 444
$ws.Range("T156").Value = 18

# --- Row 157 ---
$ws.Range("A157").Value = 11
$ws.Range("B157").Value = 'Vega Monumental Concepción'
$ws.Range("C157").Value = 'Bíobío'
$ws.Range("D157").Value = 44748
$ws.Range("E157").Value = 8
$ws.Range("F157").Value = 'Fruta'
$ws.Range("G157").Value = 100101
$ws.Range("H157").Value = 'Berries'
$ws.Range("I157").Value = 100101007
$ws.Range("J157").Value = 'Kiwi'
$ws.Range("K157").Value = 'Hayward'
$ws.Range("L157").Value = 'Segunda'
$ws.Range("M157").Value = 50
$ws.Range("N157").Value = 7000
$ws.Range("O157").Value = 7000
$ws.Range("P157").Value = 7000
$ws.Range("Q157").Value = '$/bandeja 18 kilos'
$ws.Range("R157").Value = 'Región de O''Higgins'
$ws.Range("S157").Value = 389
$ws.Range("T157").Value = 18

# --- Row 166 ---
$ws.Range("A166").Value = 11
$ws.Range("B166").Value = 'Vega Monumental Concepción'
$ws.Range("C166").Value = 'Bíobío'
$ws.Range("D166").Value = 44747
$ws.Range("E166").Value = 8
$ws.Range("F166").Value = 'Fruta'
$ws.Range("G166").Value = 100101
$ws.Range("H166").Value = 'Berries'
$ws.Range("I166").Value = 100101007
$ws.Range("J166").Value = 'Kiwi'
$ws.Range("K166").Value = 'Hayward'
$ws.Range("L166").Value = 'Primera'
$ws.Range("M166").Value = 150
$ws.Range("N166").Value = 8000
$ws.Range("O166").Value = 9000
$ws.Range("P166").Value = 8333
$ws.Range("Q166").Value = '$/bandeja 18 kilos'
$ws.Range("R166").Value = 'Región de O''Higgins'
$ws.Range("S166").Value = 463
$ws.Range("T166").Value = 18

# --- Row 167 ---
$ws.Range("A167").Value = 11
$ws.Range("B167").Value = 'Vega Monumental Concepción'
$ws.Range("C167").Value = 'Bíobío'
$ws.Range("D167").Value = 44747
$ws.Range("E167").Value = 8
$ws.Range("F167").Value = 'Fruta'
$ws.Range("G167").Value = 100101
$ws.Range("H167").Value = 'Berries'
$ws.Range("I167").Value = 100101007
$ws.Range("J167").Value = 'Kiwi'
$ws.Range("K167").Value = 'Hayward'
$ws.Range("L167").Value = 'Segunda'
$ws.Range("M167").Value = 100
$ws.Range("N167").Value = 7000
$ws.Range("O167").Value = 7000
$ws.Range("P167").Value = 7000
$ws.Range("Q167").Value = '$/bandeja 18 kilos'
$ws.Range("R167").Value = 'Región de O''Higgins'
$ws.Range("S167").Value = 389
$ws.Range("T167").Value = 18

